# Applies the "User Import Template" rework:
#  - renames the worksheet and the title string
#  - adds a "graduationYear" column (H) to the header row
#  - removes the now-unused template rows (3-23), keeping only the two header rows
#  - drops the frozen header pane
#  - resizes the columns to their new widths
#  - strips the bold/filled header styling down to a plain, non-wrapping style
#    and normalizes the header border colors
#  - widens the merged title band to A1:H1 and resets the page margins

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Sheet name / title text
# ---------------------------------------------------------------------------
$ws.Name = "Sheet 1 - User Import Template"
$ws.Range("A1").Value = "User Import Template"

# ---------------------------------------------------------------------------
# 2. Remove the sample/placeholder rows (3-23); only the title + header rows
#    (1-2) remain afterwards.
# ---------------------------------------------------------------------------
$ws.Rows("3:23").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3. Add the new "graduationYear" header in column H, matching the look of
#    the existing header cells, and extend the title merge/band to H.
# ---------------------------------------------------------------------------
$ws.Range("A1:G1").UnMerge() | Out-Null
$ws.Range("H1").Value = ""
$ws.Range("A1:H1").Merge() | Out-Null

$ws.Range("H2").Value = "graduationYear"

# ---------------------------------------------------------------------------
# 4. Remove the frozen header pane.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false

# ---------------------------------------------------------------------------
# 5. Column widths. The COM layer quantizes ColumnWidth to 1/7-character
#    steps before it is written back out as the OOXML "width" attribute
#    (width = (Round(ColumnWidth*7) + 5) / 7), so the helper below picks the
#    ColumnWidth value that reproduces each target width as closely as that
#    grid allows.
# ---------------------------------------------------------------------------
function Get-ColumnWidthFor($targetWidth) {
    $px = [Math]::Round(($targetWidth * 7) - 5)
    return $px / 7.0
}

$colTargets = @(9.35156, 10.8516, 5.5, 8.85156, 4.17188, 6.5, 12.6719, 7.35156)
for ($i = 0; $i -lt $colTargets.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = (Get-ColumnWidthFor $colTargets[$i])
}

$defaultColWidth = Get-ColumnWidthFor 8.35156
$ws.Range($ws.Columns.Item(9), $ws.Columns.Item(256)).ColumnWidth = $defaultColWidth

# ---------------------------------------------------------------------------
# 6. Header-row styling: drop the bold font + solid fill + wrap text, switch
#    the header borders to the lighter gray/charcoal combination used after
#    the cleanup.
# ---------------------------------------------------------------------------
$headerRow = $ws.Range("A2:H2")
$headerRow.Font.Bold = $false
$headerRow.Interior.ColorIndex = -4142   # xlColorIndexNone
$headerRow.Interior.Pattern = -4142      # xlPatternNone
$headerRow.WrapText = $false

foreach ($addr in @("A2", "B2", "C2", "D2", "E2", "F2", "G2", "H2")) {
    $cell = $ws.Range($addr)
    $cell.Borders.LineStyle = 1
    $cell.Borders.Item(7).Color = 0xA5A5A5   # left
    $cell.Borders.Item(8).Color = 0xA5A5A5   # top
    $cell.Borders.Item(9).Color = 0xA5A5A5   # bottom
    $cell.Borders.Item(10).Color = 0xA5A5A5  # right
}
# Column A keeps a darker accent on its right edge, column B on its left edge
$ws.Range("A2").Borders.Item(10).Color = 0x3F3F3F
$ws.Range("B2").Borders.Item(7).Color = 0x3F3F3F

# ---------------------------------------------------------------------------
# 7. Page margins (inches: 1/1/1/1 margins, 0.25 header/footer).
# ---------------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.LeftMargin = 72
$ps.RightMargin = 72
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 18
$ps.FooterMargin = 18
